# Update cryptocurrency Price (D) and Volume(1h) (E) columns with latest scraped values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.662.41'
$ws.Range("E2").Value = '  -0.62%  '
$ws.Range("D3").Value = '3.854.08'
$ws.Range("E3").Value = '  +2.92%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '599.81'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.30%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '162.05'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.07%  '
$ws.Range("D7").Value = '3.850.74'
$ws.Range("E7").Value = '  +2.89%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.529'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.42%  '
$ws.Range("E10").Value = '  -2.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.29'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.48%  '
$ws.Range("E12").Value = '  -0.45%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '36.63'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.55%  '
$ws.Range("E14").Value = '  -2.37%  '
$ws.Range("D15").Value = '4.499.73'
$ws.Range("E15").Value = '  +3.01%  '
$ws.Range("D16").Value = '3.848.62'
$ws.Range("E16").Value = '  +2.40%  '
$ws.Range("D17").Value = '68.832.08'
$ws.Range("E17").Value = '  -0.39%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.53'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.63%  '
$ws.Range("E19").Value = '  -0.56%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.04'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.77%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.26'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.16%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '483.72'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.84%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.716'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.65%  '
$ws.Range("E24").Value = '  +6.00%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.79'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.27%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.22'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.27%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.06'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.67%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.998'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.13%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.92'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.43%  '
$ws.Range("E30").Value = '  -1.18%  '
$ws.Range("D31").Value = '4.008.42'
$ws.Range("E31").Value = '  +3.02%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.83'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.80%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.99'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.65%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.35'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.33%  '
$ws.Range("D35").Value = '3.802.36'
$ws.Range("E35").Value = '  +3.40%  '
$ws.Range("E36").Value = '  -1.79%  '
$ws.Range("E37").Value = '  +1.63%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.139'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.51%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.85'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.60%  '
$ws.Range("E40").Value = '  -0.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.316'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.94%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.96'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.92%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '428.22'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.49%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '48.48'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.47%  '
$ws.Range("E45").Value = '  -0.86%  '
$ws.Range("E46").Value = '  -0.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.36'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.11%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '143.44'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.12%  '
$ws.Range("D49").Value = '2.832.63'
$ws.Range("E49").Value = '  +1.63%  '
$ws.Range("E50").Value = '  +0.99%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '25.88'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +13.22%  '
